$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (shared strings with rich-text runs) ---
# A8: "Volume 30   Number  48" -> "...49"
$ws.Range("A8").Characters(21,2).Text = "49"
# C9: "Report Covering the Week  11/27/2023  Through  12/3/2023"
#     -> "...12/4/2023  Through  12/10/2023"
# Replace the second (later) date first so the first date's character offset
# is unaffected by the length change of the replacement text.
$ws.Range("C9").Characters(48,9).Text = "12/10/2023"
$ws.Range("C9").Characters(27,10).Text = "12/4/2023"

# Row 15 (Rape)
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 17
$ws.Range("K15").Value = -15
$ws.Range("L15").Value = -19.047619047619
$ws.Range("M15").Value = 41.666666666666
$ws.Range("N15").Value = -46.875

# Row 16 (Robbery)
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 171
$ws.Range("J16").Value = 195
$ws.Range("K16").Value = -12.307692307692
$ws.Range("L16").Value = 54.054054054054
$ws.Range("M16").Value = -32.941176470588
$ws.Range("N16").Value = -84.920634920634

# Row 17 (Fel.Assault)
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -37.931034482758
$ws.Range("I17").Value = 353
$ws.Range("J17").Value = 360
$ws.Range("K17").Value = -1.944444444444
$ws.Range("L17").Value = 40.637450199203
$ws.Range("M17").Value = 54.824561403508
$ws.Range("N17").Value = 1.146131805157

# Row 18 (Burglary)
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 7.142857142857
$ws.Range("I18").Value = 133
$ws.Range("J18").Value = 149
$ws.Range("K18").Value = -10.738255033557
$ws.Range("L18").Value = 29.126213592233
$ws.Range("M18").Value = -61.449275362318
$ws.Range("N18").Value = -90.383224873463

# Row 19 (Gr.Larceny)
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 10.344827586206
$ws.Range("I19").Value = 321
$ws.Range("J19").Value = 358
$ws.Range("K19").Value = -10.335195530726
$ws.Range("L19").Value = 18.014705882352
$ws.Range("M19").Value = -5.309734513274
$ws.Range("N19").Value = -41.423357664233

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 6
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 298
$ws.Range("K20").Value = 21.138211382113
$ws.Range("L20").Value = 31.858407079646
$ws.Range("M20").Value = 8.759124087591
$ws.Range("N20").Value = -90.955993930197

# Row 21 (TOTAL)
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 9.090909090909
$ws.Range("F21").Value = 109
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = 3.809523809523
$ws.Range("I21").Value = 1300
$ws.Range("J21").Value = 1333
$ws.Range("K21").Value = -2.475618904726
$ws.Range("L21").Value = 31.445904954499
$ws.Range("M21").Value = -10.836762688614
$ws.Range("N21").Value = -80.746445497630

# Row 22 (Transit)
$ws.Range("C15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 14
$ws.Range("K22").Value = 7.692307692307
$ws.Range("L22").Value = 75
$ws.Range("M22").Value = -12.5

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -46.666666666666
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 136
$ws.Range("H24").Value = -25.735294117647
$ws.Range("I24").Value = 1301
$ws.Range("J24").Value = 1563
$ws.Range("K24").Value = -16.762635956493
$ws.Range("L24").Value = 32.215447154471
$ws.Range("M24").Value = 73.466666666666

# Row 25 (Misd.Assault)
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 28.571428571428
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 42.857142857142
$ws.Range("I25").Value = 540
$ws.Range("J25").Value = 516
$ws.Range("K25").Value = 4.651162790697
$ws.Range("L25").Value = 9.311740890688
$ws.Range("M25").Value = -8.163265306122

# Row 26 (UCR Rape*)
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("I26").Value = 36
$ws.Range("J26").Value = 36
$ws.Range("L26").Value = 9.090909090909

# Row 27 (Other Sex Crimes)
$ws.Range("C15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 51
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 2
$ws.Range("L27").Value = -26.086956521739

# Row 28 (Shooting Vic.)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4163)

# Row 29 (Shooting Inc.)
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4163)

# Row 30 (Hate Crimes)
$ws.Range("C15").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 6
$ws.Range("K30").Value = 20
$ws.Range("L30").Value = -14.285714285714
